$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(11, 2).Value = "6ad376482c35748fd2e6fc3c1f5d0554"
$ws.Cells.Item(17, 2).Value = "80f527c2ab700ca7a24fe6be2b570b8a"
$ws.Cells.Item(24, 2).Value = "9ba924a16abfcb8ef6c9e3cac8a2cd90"
$ws.Cells.Item(29, 2).Value = "2c8d30f99a54a5abe339d9c67c3b74c8"
$ws.Cells.Item(34, 2).Value = "222ef74b04f6155c4eaff1fa46b91938"
$ws.Cells.Item(121, 2).Value = "57bd65bb24119e1000f806bde286b5a6"
$ws.Cells.Item(126, 2).Value = "95a722eb6c75a8d2e7f8e464fd2caea6"
$ws.Cells.Item(133, 2).Value = "3219a7f142c0467e3e9d75510915bb73"
$ws.Cells.Item(136, 2).Value = "a64432ebb8744fa1c996a27fef29b1da"
$ws.Cells.Item(159, 2).Value = "5516c0461909e150764e8279d36584a4"
$ws.Cells.Item(162, 2).Value = "8085d16dd06d76ea05474ceeeb51071a"
$ws.Cells.Item(169, 2).Value = "6e2a5978a24279e1171507a29d14efd5"
$ws.Cells.Item(175, 2).Value = "3da1c8cefe8898bd7607a9bd90b90922"
$ws.Cells.Item(180, 2).Value = "771487aa79f983aa0733084642b626c8"
$ws.Cells.Item(183, 2).Value = "76441d79305d0d22a2eb1099c779ca26"
$ws.Cells.Item(191, 2).Value = "295819ab10107e5b676516d3e1b806e6"
$ws.Cells.Item(198, 2).Value = "307897560bcd379a5f70f3a0b70817fc"
$ws.Cells.Item(213, 2).Value = "888800f7c712b27bc876e1c1407d396b"
$ws.Cells.Item(227, 2).Value = "21b3f3ab661d6bf04cf75158df6b1f61"
$ws.Cells.Item(232, 2).Value = "5eaa4c802a4dd56ecaea734651e35d32"
$ws.Cells.Item(281, 2).Value = "17cfa0728bacabad7c7d2276ad59d422"
$ws.Cells.Item(302, 2).Value = "42808833767abeb10c32e7d7d28b5776"
$ws.Cells.Item(339, 2).Value = "bb925f9bce4146dbc18f0ef0f1387cf2"
$ws.Cells.Item(460, 2).Value = "a5dbe54c39a9069dfff780add106e62d"
$ws.Cells.Item(461, 2).Value = "868ba5fe7dbb51fb23b3548048929449"
$ws.Cells.Item(478, 2).Value = "9f0e0bf032466e2476527cbc02a4f370"
$ws.Cells.Item(480, 2).Value = "0fe57b3149dac462344231936f3e459f"
$ws.Cells.Item(500, 2).Value = "699658c5c4dee4e8bbbd60f12d5ecc22"
$ws.Cells.Item(501, 2).Value = "63fe587f91e1dc1dc5cb868231e9cf75"
$ws.Cells.Item(502, 2).Value = "55202e3ac950dec6c1cd548634aef598"
$ws.Cells.Item(506, 2).Value = "90d04f32548b90e165944f5101de222b"
$ws.Cells.Item(514, 2).Value = "d9b38b5fdf93f1853eb7fb2695d11876"
$ws.Cells.Item(515, 2).Value = "9d3a56fb1b3e81c744532ce360c2d639"
$ws.Cells.Item(517, 2).Value = "790260fc1b06d5c1e5750256043dad45"
$ws.Cells.Item(524, 2).Value = "7866d68d2e0c2a9497456576e556d779"
$ws.Cells.Item(547, 2).Value = "4faa924b5230286d8ff2c0682e53d28e"
$ws.Cells.Item(550, 2).Value = "8aa17d499f1be35da12d989f536cfad5"
$ws.Cells.Item(563, 2).Value = "d0ea26a7c1144555f02abc95e1d5b8cb"
$ws.Cells.Item(572, 2).Value = "cb5e9c0af6814178eddf0cf60e6d737b"
$ws.Cells.Item(616, 2).Value = "0a8197a280321a7f99dd9c791f024dce"
$ws.Cells.Item(627, 2).Value = "2521330e9c43a86a2061c5c26fcd442a"
$ws.Cells.Item(629, 2).Value = "6dcbe8996db1c86df41fd864e22d9d5e"
$ws.Cells.Item(649, 2).Value = "408877f30da306439518d7d711a22846"
$ws.Cells.Item(655, 2).Value = "0bec3d584237edb7911345244235a1d6"
$ws.Cells.Item(666, 2).Value = "43dc792658b0505cbbb84c25ce00acba"
$ws.Cells.Item(700, 2).Value = "efbe83cbeb5b4267c27ca6d63394791b"
$ws.Cells.Item(704, 2).Value = "0ca4ca26d3a319e633ec25db71a3c100"
$ws.Cells.Item(715, 2).Value = "acf899582af4884e0190ac6a5bf22e98"
$ws.Cells.Item(729, 2).Value = "08d1f68a176e207d91813825be6d1365"
$ws.Cells.Item(742, 2).Value = "2a1b0d7d574bd03ba3b2bd96c80e930f"
$ws.Cells.Item(819, 2).Value = "f83f8d714d1de762fed2e3e32b9de845"
$ws.Cells.Item(830, 2).Value = "8315e7a3ea6916eb2d4ec8f0540b7172"
$ws.Cells.Item(835, 2).Value = "5d16fde8cebfe5183df5a7e5d5aecfe6"
$ws.Cells.Item(854, 2).Value = "3140b1767b7d88d72d10af13a99dcc20"
$ws.Cells.Item(862, 2).Value = "661bc47e0ade0a5d4c7e5e05b5425aa6"
